# Apply "Archivos QA gateway, Archivos Excel Catalogo general" fix:
# - Rename sheet "Reactivos" -> "Catálogos"
# - Rename/resize defined name "Reactivos" -> "Catalogos", range now A4:C5 on the renamed sheet
# - Drop the "Contaq"/"Sistema" columns (C/D/E) from the header & data rows,
#   leaving Clave / Nombre / Activo as the three columns
# - Update the active cell selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Catálogos"

# Rename the workbook-level defined name and point it at the new, narrower range
$nm = $wb.Names.Item(1)
$nm.Name = "Catalogos"
$nm.RefersTo = "=Catálogos!`$A`$4:`$C`$5"

# Header row (row 3): keep Clave/Nombre, turn the old "Clave Contaq" cell into "Activo",
# blank out the old "Nombre Contaq" cell (keeping its style), and drop the old "Activo" cell (E3)
$ws.Range("C3").Value = "Activo"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Clear()

# Data row (row 4): Clave/Nombre/Activo placeholders in the first three columns,
# drop the old ClaveSistema/NombreSistema placeholder cells (D4, E4)
$ws.Range("A4").Value = "{{item.Clave}}"
$ws.Range("B4").Value = "{{item.Nombre}}"
$ws.Range("C4").Value = "{{item.Activo}}"
$ws.Range("D4").Clear()
$ws.Range("E4").Clear()

# Update the saved selection/active cell
$ws.Range("G1").Select()
